$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of "Price" values in column D are plain decimal numbers (e.g. "568.46").
# The column stores prices as text (see the other, already-dotted-thousands values such
# as "70.206.93"), so force those specific cells to Text format first. Otherwise Excel's
# smart entry would silently convert a value like "1.00" into the number 1.
$textPriceCells = @(
    "D5", "D6", "D8", "D9", "D12", "D13", "D14", "D17",
    "D19", "D21", "D23", "D25", "D26", "D27", "D31", "D32",
    "D33", "D34", "D35", "D36", "D37", "D40", "D41", "D44",
    "D45", "D46", "D47", "D48", "D49", "D50", "D51"
)
foreach ($cellRef in $textPriceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "70.206.93"
$ws.Range("E2").Value = "  +3.64%  "
$ws.Range("D3").Value = "2.461.33"
$ws.Range("E3").Value = "  +1.95%  "
$ws.Range("D5").Value = "568.46"
$ws.Range("E5").Value = "  +2.34%  "
$ws.Range("D6").Value = "168.24"
$ws.Range("E6").Value = "  +4.90%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "0.514"
$ws.Range("E8").Value = "  +0.76%  "
$ws.Range("D9").Value = "0.177"
$ws.Range("E9").Value = "  +13.76%  "
$ws.Range("D10").Value = "2.459.43"
$ws.Range("E10").Value = "  +1.86%  "
$ws.Range("E11").Value = "  -1.40%  "
$ws.Range("D12").Value = "0.337"
$ws.Range("E12").Value = "  +3.47%  "
$ws.Range("D13").Value = "4.72"
$ws.Range("E13").Value = "  -0.86%  "
$ws.Range("D14").Value = "0.0000184"
$ws.Range("E14").Value = "  +9.45%  "
$ws.Range("D15").Value = "70.069.84"
$ws.Range("E15").Value = "  +3.51%  "
$ws.Range("D16").Value = "2.910.79"
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("D17").Value = "24.20"
$ws.Range("E17").Value = "  +5.44%  "
$ws.Range("D18").Value = "2.453.60"
$ws.Range("E18").Value = "  +1.88%  "
$ws.Range("D19").Value = "10.91"
$ws.Range("E19").Value = "  +6.03%  "
$ws.Range("E20").Value = "  +5.99%  "
$ws.Range("D21").Value = "344.15"
$ws.Range("E21").Value = "  +2.92%  "
$ws.Range("D23").Value = "2.02"
$ws.Range("E23").Value = "  +8.27%  "
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").Value = "66.64"
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("D26").Value = "3.90"
$ws.Range("E26").Value = "  +8.26%  "
$ws.Range("B27").Value = "Aptos"
$ws.Range("C27").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D27").Value = "8.59"
$ws.Range("E27").Value = "  +7.11%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.581.56"
$ws.Range("E28").Value = "  +1.67%  "
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("D30").Value = "0.0₃0872"
$ws.Range("D31").Value = "7.42"
$ws.Range("E31").Value = "  +4.76%  "
$ws.Range("D32").Value = "1.26"
$ws.Range("E32").Value = "  +11.74%  "
$ws.Range("D33").Value = "456.55"
$ws.Range("E33").Value = "  +8.80%  "
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("D35").Value = "1.64"
$ws.Range("E35").Value = "  +2.66%  "
$ws.Range("D36").Value = "161.43"
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("D37").Value = "19.10"
$ws.Range("E37").Value = "  +0.86%  "
$ws.Range("E38").Value = "  +9.22%  "
$ws.Range("D40").Value = "18.24"
$ws.Range("E40").Value = "  +3.19%  "
$ws.Range("D41").Value = "0.306"
$ws.Range("E41").Value = "  +4.98%  "
$ws.Range("E42").Value = "  +7.22%  "
$ws.Range("E43").Value = "  +4.40%  "
$ws.Range("D44").Value = "37.81"
$ws.Range("E44").Value = "  +0.56%  "
$ws.Range("D45").Value = "1.10"
$ws.Range("E45").Value = "  +5.53%  "
$ws.Range("D46").Value = "2.16"
$ws.Range("E46").Value = "  +8.77%  "
$ws.Range("D47").Value = "3.42"
$ws.Range("E47").Value = "  +2.88%  "
$ws.Range("D48").Value = "133.69"
$ws.Range("E48").Value = "  +4.46%  "
$ws.Range("D49").Value = "0.0727"
$ws.Range("E49").Value = "  +2.06%  "
$ws.Range("D50").Value = "0.494"
$ws.Range("E50").Value = "  +4.23%  "
$ws.Range("D51").Value = "0.566"
$ws.Range("E51").Value = "  +2.52%  "
